# Add the newly read "pass" value into Sheet1!E2, then update the
# selection / zoom to match what Excel left behind after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data written by the reading/automation code.
$ws.Range("E2").Value = "pass"

# Leave the sheet zoomed in and selected on the new cell, matching
# the state Excel saved the workbook in.
$ws.Range("E2").Select()
$excel.ActiveWindow.Zoom = 160
